# Updates cryptos list cell values (Price column D, Volume(1h) column E)
# NumberFormat "@" + Style reset to "Normal" keeps the write as literal
# text (matching the inline-string source data, e.g. "194.40" not 194.4)
# without leaving a residual custom style behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.815.40"
Set-TextValue "E2" "  -0.26%  "

Set-TextValue "D3" "1.635.32"
Set-TextValue "E3" "  -0.15%  "

Set-TextValue "E4" "  -0.19%  "

Set-TextValue "D5" "215.51"
Set-TextValue "E5" "  +0.09%  "

Set-TextValue "E6" "  -0.85%  "

Set-TextValue "E7" "  -0.16%  "

Set-TextValue "E8" "  +0.15%  "

Set-TextValue "E9" "  -0.88%  "

Set-TextValue "D10" "19.69"
Set-TextValue "E10" "  -3.46%  "

Set-TextValue "E11" "  +1.44%  "

Set-TextValue "D12" "4.28"
Set-TextValue "E12" "  -0.13%  "

Set-TextValue "D13" "1.636.18"
Set-TextValue "E13" "  +0.00%  "

Set-TextValue "D14" "1.860.37"
Set-TextValue "E14" "  -0.17%  "

Set-TextValue "D15" "0.561"
Set-TextValue "E15" "  -0.11%  "

Set-TextValue "E16" "  +0.00%  "

Set-TextValue "D17" "62.85"
Set-TextValue "E17" "  -0.70%  "

Set-TextValue "D18" "25.817.17"
Set-TextValue "E18" "  -0.27%  "

Set-TextValue "E19" "  -0.24%  "

Set-TextValue "E20" "  +1.59%  "

Set-TextValue "D21" "194.40"
Set-TextValue "E21" "  +0.14%  "

Set-TextValue "E22" "  -0.09%  "

Set-TextValue "E23" "  +1.92%  "

Set-TextValue "E24" "  -0.15%  "

Set-TextValue "E25" "  +2.31%  "

Set-TextValue "E26" "  +3.07%  "

Set-TextValue "E27" "  -0.04%  "

Set-TextValue "E28" "  +0.83%  "

Set-TextValue "E29" "  -0.05%  "

Set-TextValue "E30" "  -0.23%  "

Set-TextValue "D31" "0.0494"
Set-TextValue "E31" "  -0.61%  "

Set-TextValue "D32" "3.36"
Set-TextValue "E32" "  +1.65%  "

Set-TextValue "E33" "  -0.25%  "

Set-TextValue "E34" "  +0.76%  "

Set-TextValue "D35" "2.38"
Set-TextValue "E35" "  -0.07%  "

Set-TextValue "E36" "  +0.07%  "

Set-TextValue "D37" "1.134.77"
Set-TextValue "E37" "  -0.19%  "

Set-TextValue "D38" "2.54"
Set-TextValue "E38" "  -1.47%  "

Set-TextValue "D39" "0.547"
Set-TextValue "E39" "  -1.74%  "

Set-TextValue "D40" "0.0156"
Set-TextValue "E40" "  -0.62%  "

Set-TextValue "E41" "  +0.24%  "

Set-TextValue "E42" "  +2.17%  "

Set-TextValue "D43" "100.59"
Set-TextValue "E43" "  +0.80%  "

Set-TextValue "E44" "  +0.69%  "

Set-TextValue "D45" "1.769.75"
Set-TextValue "E45" "  -0.39%  "

Set-TextValue "E46" "  +1.38%  "

Set-TextValue "D47" "55.29"
Set-TextValue "E47" "  -0.58%  "

Set-TextValue "D48" "0.0505"
Set-TextValue "E48" "  -0.42%  "

Set-TextValue "E49" "  -2.33%  "

Set-TextValue "E50" "  +0.28%  "

Set-TextValue "D51" "7.54"
Set-TextValue "E51" "  -2.71%  "
